$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.187.26'
$ws.Range("E2").Value = '  -1.39%  '

$ws.Range("D3").Value = '1.867.58'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("D4").Value = '''1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.57%  '

$ws.Range("D5").Value = '''236.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("D6").Value = '''1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.48%  '

$ws.Range("D7").Value = '''0.4775'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.07%  '

$ws.Range("D8").Value = '''0.2828'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.04%  '

$ws.Range("D9").Value = '''0.06515'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.10%  '

$ws.Range("D10").Value = '1.890.74'
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").Value = '''0.07368'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").Value = '''16.44'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.44%  '

$ws.Range("D13").Value = '''5.146'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.20%  '

$ws.Range("D14").Value = '''87.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.34%  '

$ws.Range("D15").Value = '''0.6480'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.33%  '

$ws.Range("D16").Value = '30.196.07'
$ws.Range("E16").Value = '  -1.20%  '

$ws.Range("D17").Value = '''13.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").Value = '''0.000007578'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.85%  '

$ws.Range("D20").Value = '2.129.64'
$ws.Range("E20").Value = '  -1.38%  '

$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").Value = '''1.008'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''5.304'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.40%  '

$ws.Range("D23").Value = '''216.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.98%  '

$ws.Range("D24").Value = '''6.106'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.95%  '

$ws.Range("D25").Value = '''9.305'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.67%  '

$ws.Range("D26").Value = '''164.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.23%  '

$ws.Range("D27").Value = '''18.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D28").Value = '''1.907'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.91%  '

$ws.Range("D29").Value = '''1.445'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("D30").Value = '''4.239'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.88%  '

$ws.Range("D31").Value = '''0.09158'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.53%  '

$ws.Range("D32").Value = '''3.963'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.14%  '

$ws.Range("D33").Value = '''0.05013'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.37%  '

$ws.Range("D34").Value = '''0.7414'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.67%  '

$ws.Range("D35").Value = '''1.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.12%  '

$ws.Range("D36").Value = '''2.693'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.76%  '

$ws.Range("D37").Value = '''0.01828'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.42%  '

$ws.Range("D38").Value = '''2.628'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.54%  '

$ws.Range("D39").Value = '''0.9060'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.86%  '

$ws.Range("D40").Value = '''2.050'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.85%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''5.923'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''106.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.31%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.005'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.16%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.4252'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.11%  '

$ws.Range("D45").Value = '''7.458'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.82%  '

$ws.Range("D46").Value = '''1.566'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.95%  '

$ws.Range("D47").Value = '''0.1306'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.88%  '

$ws.Range("D48").Value = '''64.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -10.56%  '

$ws.Range("D49").Value = '''8.912'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.64%  '

$ws.Range("D50").Value = '''34.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.61%  '

$ws.Range("D51").Value = '''0.05719'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.93%  '
